# Applies the SWC-112 -> DASP-2 (front_running) reclassification edit.
# For a number of rows, one vulnerability that was previously counted
# under the "Other" column (K) is moved to the "front_running" column (H).
# A couple of rows also gain additional counts in columns B / D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map: cell address -> new value
$changes = @{
    "H4"  = 1
    "K4"  = 0

    "D7"  = 1
    "H7"  = 1
    "K7"  = 1

    "H8"  = 1
    "K8"  = 0

    "B9"  = 2
    "D9"  = 1
    "H9"  = 1
    "K9"  = 0

    "H12" = 1
    "K12" = 0

    "H13" = 1
    "K13" = 0

    "H14" = 1
    "K14" = 0

    "H25" = 2
    "K25" = 0

    "H27" = 1
    "K27" = 2

    "H29" = 1
    "K29" = 0

    "H35" = 1
    "K35" = 0

    "B36" = 1
    "K36" = 0

    "H38" = 1
    "K38" = 0

    "H39" = 1
    "K39" = 0

    "H40" = 1
    "K40" = 0

    "H47" = 1
    "K47" = 0

    "H51" = 1
    "K51" = 1

    "H52" = 1
    "K52" = 1
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}

$wb.Save()
